# Applies the "move design change to each task" edit:
#  - Project sheet: drop the single shared DesignChangeVariation column and
#    add a new CollisionInformationExchnage column; collapse to a single task row
#  - WorkMethod sheet: QualityRate values updated 1 -> 0.7
#  - Task sheet: new per-task DesignChangeVariation column driven by a formula

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Project sheet
# ---------------------------------------------------------------------------
$wsProject = $wb.Worksheets.Item("Project")

# Remove the old shared "DesignChangeVariation" column (column C)
$wsProject.Columns.Item(3).Delete()

# Only a single task row is kept now; drop the old rows 3 and 4
$wsProject.Range("A3:A4").EntireRow.Delete()

# Add the new CollisionInformationExchnage column at the end (column G)
$wsProject.Range("G1").Value = "CollisionInformationExchnage"
$wsProject.Columns.Item(7).ColumnWidth = 24.6

# Update the remaining data row
$wsProject.Range("A2").Value = 7
$wsProject.Range("B2").Value = 30
$wsProject.Range("C2").Value = 1
$wsProject.Range("D2").Value = 1
$wsProject.Range("E2").Value = 1
$wsProject.Range("F2").Value = 0
$wsProject.Range("G2").Value = 1

$wsProject.Activate()
$wsProject.Range("A2").Select()

# ---------------------------------------------------------------------------
# WorkMethod sheet - QualityRate (column D) drops from 1 to 0.7
# ---------------------------------------------------------------------------
$wsWorkMethod = $wb.Worksheets.Item("WorkMethod")
for ($r = 2; $r -le 11; $r++) {
    $wsWorkMethod.Cells.Item($r, 4).Value = 0.7
}

$wsWorkMethod.Activate()
$wsWorkMethod.Range("D3").Select()

# ---------------------------------------------------------------------------
# Task sheet - add per-task DesignChangeVariation column (D) = 0.1 * InitialQty
# ---------------------------------------------------------------------------
$wsTask = $wb.Worksheets.Item("Task")
$wsTask.Range("D1").Value = "DesignChangeVariation"
$wsTask.Columns.Item(3).ColumnWidth = 7.96
$wsTask.Columns.Item(4).ColumnWidth = 19.35

for ($r = 2; $r -le 51; $r++) {
    $wsTask.Range("D$r").Formula = "=0.1*C$r"
}

$wsTask.Activate()
$wsTask.Range("C4").Select()

# Leave Project as the active sheet, matching the original tab selection
$wsProject.Activate()
